$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 187 ----
$ws.Cells.Item(187, 1).Value = 'SG_DT_NONE_PLSR'
$ws.Cells.Item(187, 2).Value = 'SG+DT'
$ws.Cells.Item(187, 3).Value = 'NONE'
$ws.Cells.Item(187, 4).Value = 'PLSR'
$ws.Cells.Item(187, 5).NumberFormat = "@"
$ws.Cells.Item(187, 5).Value = '0.9769'
$ws.Cells.Item(187, 6).NumberFormat = "@"
$ws.Cells.Item(187, 6).Value = '0.0189'
$ws.Cells.Item(187, 7).NumberFormat = "@"
$ws.Cells.Item(187, 7).Value = '0.9317'
$ws.Cells.Item(187, 8).NumberFormat = "@"
$ws.Cells.Item(187, 8).Value = '0.0327'
$ws.Cells.Item(187, 9).NumberFormat = "@"
$ws.Cells.Item(187, 9).Value = '3.8258'
$ws.Cells.Item(187, 10).NumberFormat = "@"
$ws.Cells.Item(187, 10).Value = '0.0253'
$ws.Cells.Item(187, 11).NumberFormat = "@"
$ws.Cells.Item(187, 11).Value = '2.9319'
$ws.Cells.Item(187, 12).NumberFormat = "@"
$ws.Cells.Item(187, 12).Value = '3.0792'
$ws.Cells.Item(187, 13).NumberFormat = "@"
$ws.Cells.Item(187, 13).Value = '3.0797'
$ws.Cells.Item(187, 14).Value = '{''method'': [''sg'', ''dt'']}'
$ws.Cells.Item(187, 15).Value = '[None]'
$ws.Cells.Item(187, 16).Value = 'n_comp=30'
$ws.Cells.Item(187, 17).Value = '[0.72 0.72 0.73 0.73 0.72 0.72 0.7  0.7  0.7  0.7  0.64 0.64 0.63 0.63
 0.64 0.64 0.62 0.62 0.62 0.62 0.43 0.43 0.46 0.34 0.34 0.31 0.31 0.72
 0.72 0.66 0.66 0.62 0.62 0.62 0.62 0.62 0.62 0.58 0.58 0.74 0.74 0.72
 0.72 0.72 0.72 0.62 0.62 0.59 0.59 0.53 0.53 0.7  0.7  0.6  0.6  0.61
 0.61 0.54 0.54 0.52 0.52 0.51 0.51 0.71 0.71 0.55 0.55 0.63 0.63 0.6
 0.6  0.46 0.46 0.73 0.73 0.71 0.71 0.67 0.67 0.59 0.59 0.1  0.1  0.63
 0.73 0.46 0.62 0.72 0.61 0.51 0.6  0.46]'
$ws.Cells.Item(187, 18).Value = '[[0.78279259]
 [0.69842338]
 [0.77171632]
 [0.74771421]
 [0.69659771]
 [0.7344295 ]
 [0.70092695]
 [0.62185358]
 [0.76869509]
 [0.7475982 ]
 [0.6417928 ]
 [0.5765341 ]
 [0.62247245]
 [0.65316219]
 [0.64984981]
 [0.57909118]
 [0.69871158]
 [0.62589845]
 [0.63547467]
 [0.67205494]
 [0.49554376]
 [0.49504514]
 [0.40759447]
 [0.32808319]
 [0.39423872]
 [0.35321917]
 [0.35256582]
 [0.67830449]
 [0.71311754]
 [0.66894115]
 [0.68642638]
 [0.61257773]
 [0.62874943]
 [0.63482372]
 [0.60969991]
 [0.56190536]
 [0.60667876]
 [0.61929425]
 [0.56381747]
 [0.72912144]
 [0.74068946]
 [0.69558299]
 [0.68020537]
 [0.69504734]
 [0.7122183 ]
 [0.60630474]
 [0.62015733]
 [0.63073669]
 [0.6118014 ]
 [0.53451807]
 [0.56434808]
 [0.69528738]
 [0.68440327]
 [0.58516574]
 [0.56700829]
 [0.61466187]
 [0.60969899]
 [0.58839034]
 [0.54496067]
 [0.54500119]
 [0.48861668]
 [0.48629985]
 [0.45308585]
 [0.67568475]
 [0.74647821]
 [0.53153692]
 [0.5398587 ]
 [0.57635481]
 [0.63898493]
 [0.64982912]
 [0.63053781]
 [0.46048744]
 [0.47305212]
 [0.72960831]
 [0.70133272]
 [0.6787117 ]
 [0.66291125]
 [0.68885847]
 [0.66661132]
 [0.54457447]
 [0.55891569]
 [0.09186618]
 [0.09854607]
 [0.62563001]
 [0.73459577]
 [0.43090607]
 [0.6148712 ]
 [0.71165987]
 [0.60550618]
 [0.50330281]
 [0.61857474]
 [0.46353455]]'

# ---- Row 188 ----
$ws.Cells.Item(188, 1).Value = 'SG_LASSO_PLSR'
$ws.Cells.Item(188, 2).Value = 'SG'
$ws.Cells.Item(188, 3).Value = 'LASSO'
$ws.Cells.Item(188, 4).Value = 'PLSR'
$ws.Cells.Item(188, 5).NumberFormat = "@"
$ws.Cells.Item(188, 5).Value = '0.9655'
$ws.Cells.Item(188, 6).NumberFormat = "@"
$ws.Cells.Item(188, 6).Value = '0.0234'
$ws.Cells.Item(188, 7).NumberFormat = "@"
$ws.Cells.Item(188, 7).Value = '0.9635'
$ws.Cells.Item(188, 8).NumberFormat = "@"
$ws.Cells.Item(188, 8).Value = '0.0239'
$ws.Cells.Item(188, 9).NumberFormat = "@"
$ws.Cells.Item(188, 9).Value = '5.2337'
$ws.Cells.Item(188, 10).NumberFormat = "@"
$ws.Cells.Item(188, 10).Value = '0.0183'
$ws.Cells.Item(188, 11).NumberFormat = "@"
$ws.Cells.Item(188, 11).Value = '1.9707'
$ws.Cells.Item(188, 12).NumberFormat = "@"
$ws.Cells.Item(188, 12).Value = '2.0315'
$ws.Cells.Item(188, 13).NumberFormat = "@"
$ws.Cells.Item(188, 13).Value = '2.0317'
$ws.Cells.Item(188, 14).Value = '{''method'': [''sg'']}'
$ws.Cells.Item(188, 15).Value = '[{''alpha'': 0.0005, ''normalize'': False}]'
$ws.Cells.Item(188, 16).Value = 'n_comp=8'
$ws.Cells.Item(188, 17).Value = '[0.72 0.72 0.73 0.73 0.72 0.72 0.7  0.7  0.7  0.7  0.64 0.64 0.63 0.63
 0.64 0.64 0.62 0.62 0.62 0.62 0.43 0.43 0.46 0.34 0.34 0.31 0.31 0.72
 0.72 0.66 0.66 0.62 0.62 0.62 0.62 0.62 0.62 0.58 0.58 0.74 0.74 0.72
 0.72 0.72 0.72 0.62 0.62 0.59 0.59 0.53 0.53 0.7  0.7  0.6  0.6  0.61
 0.61 0.54 0.54 0.52 0.52 0.51 0.51 0.71 0.71 0.55 0.55 0.63 0.63 0.6
 0.6  0.46 0.46 0.73 0.73 0.71 0.71 0.67 0.67 0.59 0.59 0.1  0.1  0.63
 0.73 0.46 0.62 0.72 0.61 0.51 0.6  0.46]'
$ws.Cells.Item(188, 18).Value = '[[0.71545736]
 [0.70898687]
 [0.74077287]
 [0.72008441]
 [0.69949072]
 [0.73854821]
 [0.69229279]
 [0.69514524]
 [0.71828706]
 [0.71775247]
 [0.67128003]
 [0.62768502]
 [0.61091465]
 [0.60189326]
 [0.65078522]
 [0.63781563]
 [0.65129913]
 [0.64469616]
 [0.61861766]
 [0.66211137]
 [0.48635311]
 [0.49359267]
 [0.42095387]
 [0.37174963]
 [0.37633605]
 [0.31121659]
 [0.31640141]
 [0.7130768 ]
 [0.7124537 ]
 [0.63770392]
 [0.65825413]
 [0.62111706]
 [0.60131418]
 [0.64177014]
 [0.64355614]
 [0.6147127 ]
 [0.64711607]
 [0.61590529]
 [0.57754952]
 [0.73428934]
 [0.72722128]
 [0.7246355 ]
 [0.71010161]
 [0.70591437]
 [0.72931741]
 [0.60175595]
 [0.59345173]
 [0.62880794]
 [0.59716931]
 [0.54086073]
 [0.55357594]
 [0.69979692]
 [0.69897663]
 [0.59889998]
 [0.58036678]
 [0.59772809]
 [0.59780051]
 [0.59820588]
 [0.54492776]
 [0.55266434]
 [0.51469309]
 [0.50931836]
 [0.49578793]
 [0.69198045]
 [0.70600617]
 [0.54379606]
 [0.55053953]
 [0.58328099]
 [0.56962935]
 [0.63651691]
 [0.63009846]
 [0.49714545]
 [0.48320998]
 [0.71550712]
 [0.71636731]
 [0.70156642]
 [0.69392275]
 [0.64454993]
 [0.68188344]
 [0.55582664]
 [0.57515395]
 [0.09679976]
 [0.10599744]
 [0.61423429]
 [0.7130312 ]
 [0.40046609]
 [0.62431629]
 [0.72436868]
 [0.59261627]
 [0.50844822]
 [0.64221464]
 [0.49347521]]'

# ---- Row 189 ----
$ws.Cells.Item(189, 1).Value = 'SG_LASSO_PLSR'
$ws.Cells.Item(189, 2).Value = 'SG'
$ws.Cells.Item(189, 3).Value = 'LASSO'
$ws.Cells.Item(189, 4).Value = 'PLSR'
$ws.Cells.Item(189, 5).NumberFormat = "@"
$ws.Cells.Item(189, 5).Value = '0.9660'
$ws.Cells.Item(189, 6).NumberFormat = "@"
$ws.Cells.Item(189, 6).Value = '0.0232'
$ws.Cells.Item(189, 7).NumberFormat = "@"
$ws.Cells.Item(189, 7).Value = '0.9641'
$ws.Cells.Item(189, 8).NumberFormat = "@"
$ws.Cells.Item(189, 8).Value = '0.0237'
$ws.Cells.Item(189, 9).NumberFormat = "@"
$ws.Cells.Item(189, 9).Value = '5.2809'
$ws.Cells.Item(189, 10).NumberFormat = "@"
$ws.Cells.Item(189, 10).Value = '0.0184'
$ws.Cells.Item(189, 11).NumberFormat = "@"
$ws.Cells.Item(189, 11).Value = '1.9790'
$ws.Cells.Item(189, 12).NumberFormat = "@"
$ws.Cells.Item(189, 12).Value = '2.0333'
$ws.Cells.Item(189, 13).NumberFormat = "@"
$ws.Cells.Item(189, 13).Value = '2.0326'
$ws.Cells.Item(189, 14).Value = '{''method'': [''sg'']}'
$ws.Cells.Item(189, 15).Value = '[{''alpha'': 0.0005, ''normalize'': False}]'
$ws.Cells.Item(189, 16).Value = 'n_comp=11'
$ws.Cells.Item(189, 17).Value = '[0.72 0.72 0.73 0.73 0.72 0.72 0.7  0.7  0.7  0.7  0.64 0.64 0.63 0.63
 0.64 0.64 0.62 0.62 0.62 0.62 0.43 0.43 0.46 0.34 0.34 0.31 0.31 0.72
 0.72 0.66 0.66 0.62 0.62 0.62 0.62 0.62 0.62 0.58 0.58 0.74 0.74 0.72
 0.72 0.72 0.72 0.62 0.62 0.59 0.59 0.53 0.53 0.7  0.7  0.6  0.6  0.61
 0.61 0.54 0.54 0.52 0.52 0.51 0.51 0.71 0.71 0.55 0.55 0.63 0.63 0.6
 0.6  0.46 0.46 0.73 0.73 0.71 0.71 0.67 0.67 0.59 0.59 0.1  0.1  0.63
 0.73 0.46 0.62 0.72 0.61 0.51 0.6  0.46]'
$ws.Cells.Item(189, 18).Value = '[[0.72721184]
 [0.71614959]
 [0.75266576]
 [0.72391554]
 [0.70188385]
 [0.73093248]
 [0.68812234]
 [0.69074794]
 [0.72783293]
 [0.71952688]
 [0.66740802]
 [0.62576739]
 [0.60716775]
 [0.60393128]
 [0.6526802 ]
 [0.6335325 ]
 [0.65931682]
 [0.63972494]
 [0.61601348]
 [0.6542499 ]
 [0.47369218]
 [0.4893124 ]
 [0.41309739]
 [0.36684198]
 [0.37484394]
 [0.31548577]
 [0.31769943]
 [0.7082337 ]
 [0.71239181]
 [0.63346153]
 [0.65072281]
 [0.62075321]
 [0.60262732]
 [0.63618555]
 [0.64392215]
 [0.61663066]
 [0.64195517]
 [0.61281587]
 [0.579456  ]
 [0.73634388]
 [0.72745691]
 [0.72223101]
 [0.70815995]
 [0.70648097]
 [0.73420815]
 [0.60105411]
 [0.59722253]
 [0.62100803]
 [0.59317895]
 [0.54696739]
 [0.56255577]
 [0.69701529]
 [0.70059078]
 [0.58926203]
 [0.57780479]
 [0.60282081]
 [0.60021997]
 [0.60348741]
 [0.53617512]
 [0.55483658]
 [0.51709408]
 [0.50835917]
 [0.4944148 ]
 [0.69296301]
 [0.70979381]
 [0.54959373]
 [0.55298178]
 [0.58516652]
 [0.56623828]
 [0.63512301]
 [0.6298464 ]
 [0.49580453]
 [0.48439417]
 [0.71894192]
 [0.71583956]
 [0.70232787]
 [0.6907815 ]
 [0.64637198]
 [0.68852332]
 [0.55590553]
 [0.57427208]
 [0.09449377]
 [0.10618872]
 [0.61365285]
 [0.71800074]
 [0.40024276]
 [0.62368713]
 [0.72220802]
 [0.5962827 ]
 [0.50684454]
 [0.64047157]
 [0.49094882]]'

Write-Host "Added rows 187-189"
